$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.173.72'
$ws.Range("E2").Value = '  -1.81%  '
$ws.Range("D3").Value = '2.910.32'
$ws.Range("E3").Value = '  -0.75%  '
$ws.Range("E4").Value = '  -0.06%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '348.82'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.07%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '105.21'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -6.64%  '
$ws.Range("E7").Value = '  -0.99%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -3.45%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '37.44'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -5.32%  '
$ws.Range("E11").Value = '  +1.48%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.0843'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -4.46%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '18.80'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -5.99%  '
$ws.Range("D14").Value = '3.365.90'
$ws.Range("E14").Value = '  -0.74%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '7.55'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -2.63%  '
$ws.Range("D16").Value = '2.919.92'
$ws.Range("E16").Value = '  -0.64%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.951'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -3.24%  '
$ws.Range("D18").Value = '51.118.27'
$ws.Range("E18").Value = '  -2.04%  '
$ws.Range("E19").Value = '  +3.30%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '7.36'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -3.27%  '
$ws.Range("E21").Value = '  -6.21%  '
$ws.Range("D22").Value = '0.0₃0958'
$ws.Range("E22").Value = '  -2.41%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '68.51'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -3.59%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '258.87'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -3.71%  '
$ws.Range("E25").Value = '  -3.62%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.171'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -4.34%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '26.20'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -3.02%  '
$ws.Range("E28").Value = '  +0.08%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '7.28'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +4.56%  '
$ws.Range("E30").Value = '  -0.50%  '
$ws.Range("E31").Value = '  -4.72%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '6.12'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +1.22%  '
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '35.31'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -4.90%  '
$ws.Range("B34").Value = 'Toncoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.13'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -5.84%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '50.22'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -5.28%  '
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("E37").Value = '  -6.71%  '
$ws.Range("E38").Value = '  -7.81%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '17.50'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -6.05%  '
$ws.Range("E40").Value = '  -6.38%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '2.60'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -3.82%  '
$ws.Range("E42").Value = '  -2.13%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '22.24'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -3.37%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '119.76'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +8.68%  '
$ws.Range("E45").Value = '  -2.63%  '
$ws.Range("D46").Value = '2.086.92'
$ws.Range("E46").Value = '  -4.78%  '
$ws.Range("E47").Value = '  -6.77%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.27'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -9.61%  '
$ws.Range("E49").Value = '  -4.13%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0329'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -4.92%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.881'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -7.84%  '
